$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.994.57"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.633.21"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'214.02"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").Value = "'18.45"
$ws.Range("E10").Value = "  -6.06%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "1.860.96"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "'4.20"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "1.633.97"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "26.001.12"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "0.0₃0744"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "'61.79"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D20").Value = "'190.02"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "'9.57"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "'6.13"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'143.16"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "'6.76"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").Value = "'3.15"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "1.135.53"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("D38").Value = "'2.42"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "'0.525"
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "'98.55"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'0.783"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").Value = "'5.26"
$ws.Range("E43").Value = "  -4.66%  "
$ws.Range("D44").Value = "1.771.34"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "'55.13"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").Value = "'0.0524"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'7.53"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = "  +0.18%  "
